$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.942.82"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.821.46"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'355.55"
$ws.Range("E5").Value = "  +3.73%  "
$ws.Range("E6").Value = "  -2.99%  "
$ws.Range("E7").Value = "  +2.69%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.602"
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("D10").Value = "'40.99"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("D11").Value = "'0.0857"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "'0.131"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D14").Value = "'7.76"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "3.265.64"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "2.837.06"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").Value = "'0.924"
$ws.Range("E17").Value = "  +5.19%  "
$ws.Range("D18").Value = "51.872.90"
$ws.Range("D19").Value = "'7.55"
$ws.Range("E19").Value = "  +7.44%  "
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").Value = "'13.42"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("D23").Value = "'70.02"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "'268.34"
$ws.Range("E24").Value = "  -3.15%  "
$ws.Range("D25").Value = "'2.79"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").Value = "'27.06"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "'10.32"
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").Value = "'0.0475"
$ws.Range("E30").Value = "  +23.82%  "
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").Value = "'52.59"
$ws.Range("E32").Value = "  +4.66%  "
$ws.Range("D33").Value = "'34.68"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").Value = "'5.90"
$ws.Range("E34").Value = "  +3.30%  "
$ws.Range("E35").Value = "  +8.50%  "
$ws.Range("D36").Value = "'0.0845"
$ws.Range("E36").Value = "  +3.12%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'3.31"
$ws.Range("E38").Value = "  +2.52%  "
$ws.Range("E39").Value = "  -2.67%  "
$ws.Range("D40").Value = "'18.49"
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").Value = "'2.56"
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("D43").Value = "'23.25"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").Value = "'124.31"
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("E45").Value = "  -3.37%  "
$ws.Range("D46").Value = "2.094.95"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("D49").Value = "'6.00"
$ws.Range("E49").Value = "  +7.85%  "
$ws.Range("D50").Value = "'0.969"
$ws.Range("E50").Value = "  +7.85%  "
$ws.Range("D51").Value = "'9.09"
$ws.Range("E51").Value = "  +2.76%  "
